$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Ativação:" date from 01/01/2015 to 01/01/2021
$ws.Range("B8").Value = "01/01/2021"
$ws.Range("C8").Value = "01/01/2021"

# Insert a new row at 14 (shifts current rows 14-22 down to 15-23) and
# fill in the new co-instructor entry on columns B and C only.
$ws.Rows("14:14").Insert()
$ws.Range("B14").Value = "11079086 - Herlandí de Souza Andrade"
$ws.Range("C14").Value = "11079086 - Herlandí de Souza Andrade"

# Update "Método:" text (now on row 20 after the insertion)
$ws.Range("B20").Value = "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras."
$ws.Range("C20").Value = "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras."

# Update "Critério:" text (now on row 21 after the insertion)
$ws.Range("B21").Value = "Média Aritmética dos Projetos, Trabalhos e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas."
$ws.Range("C21").Value = "Média Aritmética dos Projetos, Trabalhos e Exercícios realizados no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude) desenvolvidas."

# Update "Norma de recuperação:" text (now on row 22 after the insertion)
$ws.Range("B22").Value = "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação."
$ws.Range("C22").Value = "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação."
